# "added kg to pop_stats" - Test Cases List.xlsx
# Updates a handful of test-case rows on Sheet1: marks the "seq ont all
# transcripts" rows as "just regtest", fills in the pubmed articles/search
# term rows (now Done / pop_stats... wait, pubmed), and most importantly
# marks the "1000 genomes freq" row (kg) as Done / pop_stats with a new
# comment about failures in NC regions. Also tidies a stray comment and a
# couple of priority numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 - protein seq change: drop the stale "needs upgrading" comment
$ws.Range("E19").ClearContents()

# Rows 20/21 - seq ont transcript / seq ont all transcripts: simplify comment
$ws.Range("E20").Value = "just regtest"
$ws.Range("E21").Value = "just regtest"

# Row 38 - target: clarify the CRAVAT comment spacing
$ws.Range("E38").Value = "CRAVAT db not matching any available dbs.     xls files in testing folder"

# Row 44 - pubmed articles: now done, parsed via pubmed module
$ws.Range("C44").Value = "yes"
$ws.Range("D44").Value = "pubmed"
$ws.Range("E44").Value = "just regtest"

# Row 45 - pubmed search term: same treatment, priority bumped to 10
$ws.Range("A45").Value = 10
$ws.Range("C45").Value = "yes"
$ws.Range("D45").Value = "pubmed"
$ws.Range("E45").Value = "just regtest"

# Row 46 - 1000 genomes freq: added to pop_stats, priority bumped to 10
$ws.Range("A46").Value = 10
$ws.Range("C46").Value = "yes"
$ws.Range("D46").Value = "pop_stats"
$ws.Range("E46").Value = "failures in NC regions"

# Row 58 - in tcga mutation cluster: priority corrected to 1
$ws.Range("A58").Value = 1

# Restore the scroll position / selection as left by the author
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A47").Select()
